$d = $word.ActiveDocument

$replacements = @(
    @{Old = "959÷7="; New = "563÷7="},
    @{Old = "764÷5="; New = "378÷3="},
    @{Old = "408÷4="; New = "297÷3="},
    @{Old = "363÷2="; New = "402÷6="},
    @{Old = "816÷2="; New = "480÷6="},
    @{Old = "573÷6="; New = "526÷2="},
    @{Old = "957÷9="; New = "739÷4="},
    @{Old = "697÷2="; New = "848÷8="},
    @{Old = "578÷2="; New = "425÷5="},
    @{Old = "984÷8="; New = "110÷6="},
    @{Old = "782÷8="; New = "321÷6="},
    @{Old = "893÷5="; New = "830÷8="},
    @{Old = "260÷7="; New = "591÷9="},
    @{Old = "330÷6="; New = "110÷4="},
    @{Old = "565÷4="; New = "684÷3="},
    @{Old = "796÷9="; New = "236÷6="},
    @{Old = "857÷7="; New = "513÷8="},
    @{Old = "265÷3="; New = "911÷5="},
    @{Old = "759÷3="; New = "885÷2="},
    @{Old = "694÷9="; New = "424÷8="},
    @{Old = "377÷8="; New = "251÷2="},
    @{Old = "797÷2="; New = "792÷7="},
    @{Old = "464÷7="; New = "859÷8="},
    @{Old = "347÷2="; New = "174÷5="},
    @{Old = "940÷7="; New = "578÷8="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.New, 2)
}
